$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'44.002.56"
$ws.Range("E2").Value = "  +4.28%  "
$ws.Range("D3").Value = "'2.252.58"
$ws.Range("E3").Value = "  +1.69%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'229.71"
$ws.Range("E5").Value = "  -0.32%  "
$ws.Range("D6").Value = "'0.632"
$ws.Range("E6").Value = "  +2.18%  "
$ws.Range("D7").Value = "'63.27"
$ws.Range("E7").Value = "  +4.50%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "'0.442"
$ws.Range("E9").Value = "  +10.02%  "
$ws.Range("D10").Value = "'0.101"
$ws.Range("E10").Value = "  +12.32%  "
$ws.Range("D11").Value = "'57.17"
$ws.Range("E11").Value = "  -0.50%  "
$ws.Range("D12").Value = "'25.91"
$ws.Range("E12").Value = "  +17.02%  "
$ws.Range("E13").Value = "  +1.47%  "
$ws.Range("D14").Value = "'2.594.63"
$ws.Range("E14").Value = "  +1.93%  "
$ws.Range("D15").Value = "'15.59"
$ws.Range("E15").Value = "  +0.65%  "
$ws.Range("D16").Value = "'6.13"
$ws.Range("E16").Value = "  +10.28%  "
$ws.Range("E17").Value = "  +6.39%  "
$ws.Range("D18").Value = "'2.250.94"
$ws.Range("E18").Value = "  +0.64%  "
$ws.Range("D19").Value = "'43.892.09"
$ws.Range("E19").Value = "  +4.05%  "
$ws.Range("E20").Value = "  +6.86%  "
$ws.Range("D21").Value = "'73.01"
$ws.Range("E21").Value = "  +1.45%  "
$ws.Range("D22").Value = "'6.01"
$ws.Range("E22").Value = "  -2.58%  "
$ws.Range("D23").Value = "'251.67"
$ws.Range("E23").Value = "  +3.15%  "
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("E25").Value = "  +1.32%  "
$ws.Range("E26").Value = "  -2.37%  "
$ws.Range("D27").Value = "'3.36"
$ws.Range("E27").Value = "  +26.30%  "
$ws.Range("E28").Value = "  +3.13%  "
$ws.Range("D29").Value = "'172.04"
$ws.Range("E29").Value = "  +1.51%  "
$ws.Range("D30").Value = "'20.68"
$ws.Range("E30").Value = "  +1.80%  "
$ws.Range("E31").Value = "  -2.73%  "
$ws.Range("D32").Value = "'1.37"
$ws.Range("E32").Value = "  -4.97%  "
$ws.Range("D33").Value = "'0.123"
$ws.Range("E33").Value = "  +2.50%  "
$ws.Range("D34").Value = "'0.0682"
$ws.Range("E34").Value = "  +4.85%  "
$ws.Range("D35").Value = "'4.72"
$ws.Range("E35").Value = "  +2.30%  "
$ws.Range("D36").Value = "'4.83"
$ws.Range("E36").Value = "  -3.08%  "
$ws.Range("D37").Value = "'3.79"
$ws.Range("E37").Value = "  +7.32%  "
$ws.Range("D38").Value = "'6.49"
$ws.Range("E38").Value = "  +2.35%  "
$ws.Range("D39").Value = "'2.29"
$ws.Range("E39").Value = "  -1.72%  "
$ws.Range("E40").Value = "  +2.26%  "
$ws.Range("E41").Value = "  +0.09%  "
$ws.Range("D42").Value = "'17.29"
$ws.Range("E42").Value = "  +9.10%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'8.23"
$ws.Range("E43").Value = "  -3.78%  "
$ws.Range("B44").Value = "Cronos"
$ws.Range("C44").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D44").Value = "'0.0966"
$ws.Range("E44").Value = "  +0.97%  "
$ws.Range("D45").Value = "'0.000213"
$ws.Range("E45").Value = "  -7.97%  "
$ws.Range("D46").Value = "'97.04"
$ws.Range("E46").Value = "  +0.44%  "
$ws.Range("E47").Value = "  -0.45%  "
$ws.Range("D48").Value = "'4.34"
$ws.Range("E48").Value = "  -0.63%  "
$ws.Range("D49").Value = "'1.435.60"
$ws.Range("E49").Value = "  -1.38%  "
$ws.Range("E50").Value = "  +2.95%  "
$ws.Range("E51").Value = "  -0.06%  "